$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.919.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.995.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.604"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.13%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.20"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.17%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.24%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.289.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.00"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.96"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.733"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.03%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -5.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.988.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.862.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.64%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0807"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "221.72"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -8.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.64%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.75%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.64%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.32%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.51%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.58%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.21"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.459.34"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "93.99"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.93%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.19"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.43%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.73%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.82%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.179.04"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.83"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.34%  "
